$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-11-14 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-15 Wednesday", 2)

$d.Content.Find.Execute("75×71=5325", $true, $false, $false, $false, $false, $true, 1, $false, "65×88=5720", 2)
$d.Content.Find.Execute("17×92=1564", $true, $false, $false, $false, $false, $true, 1, $false, "59×33=1947", 2)
$d.Content.Find.Execute("22×46=1012", $true, $false, $false, $false, $false, $true, 1, $false, "68×11=748", 2)
$d.Content.Find.Execute("64×76=4864", $true, $false, $false, $false, $false, $true, 1, $false, "88×86=7568", 2)
$d.Content.Find.Execute("11×38=418", $true, $false, $false, $false, $false, $true, 1, $false, "88×39=3432", 2)
$d.Content.Find.Execute("49×59=2891", $true, $false, $false, $false, $false, $true, 1, $false, "64×84=5376", 2)
$d.Content.Find.Execute("23×72=1656", $true, $false, $false, $false, $false, $true, 1, $false, "62×52=3224", 2)
$d.Content.Find.Execute("48×23=1104", $true, $false, $false, $false, $false, $true, 1, $false, "39×50=1950", 2)
$d.Content.Find.Execute("98×20=1960", $true, $false, $false, $false, $false, $true, 1, $false, "84×49=4116", 2)
$d.Content.Find.Execute("80×27=2160", $true, $false, $false, $false, $false, $true, 1, $false, "32×87=2784", 2)
$d.Content.Find.Execute("21×20=420", $true, $false, $false, $false, $false, $true, 1, $false, "32×51=1632", 2)
$d.Content.Find.Execute("15×77=1155", $true, $false, $false, $false, $false, $true, 1, $false, "70×50=3500", 2)
$d.Content.Find.Execute("74×94=6956", $true, $false, $false, $false, $false, $true, 1, $false, "52×55=2860", 2)
$d.Content.Find.Execute("24×69=1656", $true, $false, $false, $false, $false, $true, 1, $false, "25×42=1050", 2)
$d.Content.Find.Execute("66×47=3102", $true, $false, $false, $false, $false, $true, 1, $false, "47×27=1269", 2)
$d.Content.Find.Execute("60×78=4680", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=4524", 2)
$d.Content.Find.Execute("64×45=2880", $true, $false, $false, $false, $false, $true, 1, $false, "63×35=2205", 2)
$d.Content.Find.Execute("76×35=2660", $true, $false, $false, $false, $false, $true, 1, $false, "81×22=1782", 2)
$d.Content.Find.Execute("20×28=560", $true, $false, $false, $false, $false, $true, 1, $false, "45×34=1530", 2)
$d.Content.Find.Execute("41×99=4059", $true, $false, $false, $false, $false, $true, 1, $false, "40×46=1840", 2)
$d.Content.Find.Execute("93×48=4464", $true, $false, $false, $false, $false, $true, 1, $false, "24×15=360", 2)
$d.Content.Find.Execute("97×40=3880", $true, $false, $false, $false, $false, $true, 1, $false, "43×47=2021", 2)
$d.Content.Find.Execute("35×31=1085", $true, $false, $false, $false, $false, $true, 1, $false, "17×93=1581", 2)
$d.Content.Find.Execute("99×64=6336", $true, $false, $false, $false, $false, $true, 1, $false, "72×65=4680", 2)
$d.Content.Find.Execute("19×11=209", $true, $false, $false, $false, $false, $true, 1, $false, "31×83=2573", 2)
